# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 08:52"

# Estados Unidos (row 4): update Casos activos / Recuperados
$ws.Range("D4").Value = 14828
$ws.Range("E4").Value = 288355

# Chequia (row 27): update Casos totales / Nuevos casos / Muertes hoy / Muertes
$ws.Range("B27").Value = 4475
$ws.Range("C27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 62

# Reorder Mexico / Finlandia (rows 42-43) and refresh their figures.
# Row 42 becomes Finlandia, row 43 becomes Mexico.
$ws.Range("A42").Value = "Finlandia"
$ws.Range("B42").Value = 1927
$ws.Range("C42").Value = 45
$ws.Range("D42").Value = 300
$ws.Range("E42").Value = 1602
$ws.Range("F42").Value = 73
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 25

$ws.Range("A43").Value = "Mexico"
$ws.Range("B43").Value = 1890
$ws.Range("C43").Value = 202
$ws.Range("D43").Value = 633
$ws.Range("E43").Value = 1178
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 19
$ws.Range("H43").Value = 79

# Bosnia y Herzegovina (row 72): update Casos totales / Nuevos casos / Recuperados
$ws.Range("B72").Value = 626
$ws.Range("C72").Value = 2
$ws.Range("E72").Value = 575

# Reorder Sri Lanka / Georgia (rows 109-110) and refresh their figures.
# Row 109 becomes Georgia, row 110 becomes Sri Lanka.
$ws.Range("A109").Value = "Georgia"
$ws.Range("B109").Value = 170
$ws.Range("C109").Value = 8
$ws.Range("D109").Value = 36
$ws.Range("E109").Value = 133
$ws.Range("F109").Value = 6
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 1

$ws.Range("A110").Value = "Sri Lanka"
$ws.Range("B110").Value = 166
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 27
$ws.Range("E110").Value = 134
$ws.Range("F110").Value = 5
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 5
